# Updated with new release version 1.5.0
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "ContentType" column (C) for the REST API
# example rows with the value "api".
$ws.Range("C2").Value = "api"
$ws.Range("C3").Value = "api"
$ws.Range("C4").Value = "api"
$ws.Range("C7").Value = "api"

# Add a hyperlink on the URL cell E2, pointing at the same address already
# shown as its text (same target used by the other "persons/bgates" rows).
$e2 = $ws.Range("E2")
$ws.Hyperlinks.Add($e2, "https://live.virtualandemo.com/api/persons/bgates")

# Adding the hyperlink makes Excel re-apply the built-in "Hyperlink" cell
# style to E2. The cell was already formatted that way (column E uses the
# Hyperlink style throughout), so restore its original formatting by
# copying it from a sibling cell that already has that exact formatting.
$ws.Range("E4").Copy()
$e2.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Update the visible selection / scroll position saved with the sheet.
$ws.Range("E16").Select() | Out-Null
